# Update row 20 (2025Q2) metrics in metricas_recorrencia_trimestral
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C20").Value = 347
$ws.Range("D20").Value = 269
$ws.Range("E20").Value = 78
$ws.Range("F20").Value = 83.28173374613003
